# "On a finalement le sql de la BD !!!"
#
# Mark the rows about "faire les relations et identifier les types" /
# "convertir le shema" (rows 25-27) and "generer les cles etrangeres" /
# "generer les regles referentielles" (rows 30-32) as done by putting an
# "X" in column C, exactly like all the other completed steps already on
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$doneRows = @(25, 26, 27, 30, 31, 32)
foreach ($r in $doneRows) {
    $ws.Range("C" + $r).Value = "X"
}

# Scroll the frozen view down so the newly-finished rows are in view, and
# leave the selection on the last touched cell (C33), matching where the
# author's cursor ended up after ticking these boxes off.
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$ws.Range("C33").Select()
